$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextCell 2 4 '59.713.04'
Set-TextCell 2 5 '  -1.91%  '
Set-TextCell 3 4 '2.603.60'
Set-TextCell 3 5 '  +0.39%  '
Set-TextCell 4 5 '  -0.01%  '
Set-TextCell 5 4 '514.18'
Set-TextCell 5 5 '  -1.68%  '
Set-TextCell 6 4 '146.13'
Set-TextCell 6 5 '  -5.43%  '
Set-TextCell 7 4 '1.00'
Set-TextCell 7 5 '  +0.10%  '
Set-TextCell 8 4 '0.562'
Set-TextCell 8 5 '  -4.93%  '
Set-TextCell 9 4 '2.603.37'
Set-TextCell 9 5 '  +0.15%  '
Set-TextCell 10 4 '6.23'
Set-TextCell 10 5 '  -7.13%  '
Set-TextCell 12 5 '  -3.42%  '
Set-TextCell 13 5 '  -0.89%  '
Set-TextCell 14 4 '3.057.90'
Set-TextCell 14 5 '  +0.28%  '
Set-TextCell 15 4 '59.705.47'
Set-TextCell 15 5 '  -1.95%  '
Set-TextCell 16 4 '20.89'
Set-TextCell 16 5 '  -3.44%  '
Set-TextCell 17 5 '  -3.56%  '
Set-TextCell 18 4 '2.608.37'
Set-TextCell 18 5 '  +0.32%  '
Set-TextCell 19 5 '  -3.59%  '
Set-TextCell 20 4 '337.66'
Set-TextCell 20 5 '  -4.25%  '
Set-TextCell 21 4 '10.24'
Set-TextCell 21 5 '  -3.22%  '
Set-TextCell 22 4 '6.02'
Set-TextCell 22 5 '  -3.28%  '
Set-TextCell 23 5 '  -0.05%  '
Set-TextCell 24 4 '60.60'
Set-TextCell 24 5 '  -0.61%  '
Set-TextCell 25 5 '  -3.11%  '
Set-TextCell 26 4 '1.00'
Set-TextCell 26 5 '  +0.13%  '
Set-TextCell 27 5 '  -5.22%  '
Set-TextCell 28 4 '0.0₃0789'
Set-TextCell 28 5 '  -6.56%  '
Set-TextCell 29 5 '  -5.91%  '
Set-TextCell 30 5 '  +0.02%  '
Set-TextCell 31 4 '1.56'
Set-TextCell 31 5 '  -2.39%  '
Set-TextCell 32 4 '5.90'
Set-TextCell 32 5 '  -6.60%  '
Set-TextCell 33 4 '18.73'
Set-TextCell 33 5 '  -3.17%  '
Set-TextCell 34 4 '149.91'
Set-TextCell 34 5 '  +1.25%  '
Set-TextCell 35 4 '3.88'
Set-TextCell 35 5 '  -7.06%  '
Set-TextCell 36 4 '0.896'
Set-TextCell 36 5 '  -3.63%  '
Set-TextCell 37 5 '  -6.70%  '
Set-TextCell 38 4 '36.56'
Set-TextCell 38 5 '  +0.23%  '
Set-TextCell 39 4 '0.840'
Set-TextCell 39 5 '  -0.98%  '
Set-TextCell 40 5 '  -6.38%  '
Set-TextCell 41 4 '3.56'
Set-TextCell 41 5 '  -6.00%  '
Set-TextCell 42 4 '284.41'
Set-TextCell 42 5 '  -1.18%  '
Set-TextCell 43 4 '0.618'
Set-TextCell 43 5 '  -0.65%  '
Set-TextCell 44 5 '  +0.15%  '
Set-TextCell 45 5 '  -2.23%  '
Set-TextCell 46 4 '0.0537'
Set-TextCell 46 5 '  -4.02%  '
Set-TextCell 47 4 '18.89'
Set-TextCell 47 5 '  -3.95%  '
Set-TextCell 48 5 '  +0.51%  '
Set-TextCell 49 5 '  -3.37%  '
Set-TextCell 50 4 '4.55'
Set-TextCell 50 5 '  -6.82%  '
Set-TextCell 51 4 '1.929.92'
Set-TextCell 51 5 '  -1.10%  '
